$d = $word.ActiveDocument

# Locate the paragraph ending "...regardless of the data." which is the
# anchor point: five new paragraphs about the fractional penalty are
# inserted immediately after it (and before the blank paragraphs that
# already follow it).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*We had a choice for how to limit the number of splits performed*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find anchor paragraph"
}

# Build a zero-width range positioned just before the anchor paragraph's
# end-of-paragraph mark (NOT Collapse(0), which seats the caret against
# the paragraph-mark boundary and causes InsertXML to swallow the
# adjoining paragraph). Anchoring one character earlier keeps the insert
# squarely inside the anchor paragraph's story, so InsertXML splits the
# story there and leaves every neighboring paragraph untouched.
$insertAt = $target.Range.End - 1
$r = $d.Range($insertAt, $insertAt)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:after="160"/></w:pPr><w:r><w:t xml:space="preserve">But we were still faced with a choice: which type of penalty to use? We could set a fixed penalty, such that the SSE would need to be reduced by a certain, fixed amount in order to justify a split. Since we would allow this fixed number to be tuned for each tree, we could achieve a reasonable fit this way. Or, we could implement a fractional penalty, which is even more flexible. It would require that each split reduce the SSE by at least a certain fraction, probably in the 0.05-0.15 range. If the split failed to reduce the SSE by that fraction, the split would not </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>occur</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> and the node would be marked as a terminal node.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="160"/></w:pPr><w:r><w:t>Given the increased flexibility that the fractional penalty provided, this seemed like the best choice.</w:t></w:r><w:r><w:t xml:space="preserve"> In a regression tree, it is inevitable that some nodes will contain observations that are very similar, while others contain observations that are only somewhat similar. For example, at the very first split, when all the data is under consideration, it is possible that there really are two clearly distinct groups in the data. The first group represents about half of the data and consists of nearly identical observations. The second group, which likewise represents half of the data, is remarkable for the relative dissimilarity of its observations. Assume that the split is performed correctly, so that the two groups are situated in two different nodes. The first node may very well require no further splits. The second node, on the other hand, consists of observations which have not much in common. It was formed merely by default, the leftover result of the split which so accurately segregated the first group. Here we have a case where </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t>some nodes contain observations which are very similar, while some nodes contain observations which are very different.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="160"/></w:pPr><w:r><w:t xml:space="preserve">If a node contains observations which are very different, we would want those observations to be further split out, until we had nodes which contained relatively similar observations. In the case of a fixed penalty, this may not occur. A large node may contain relatively similar observations, but because of its sheer size, it reduces the SSE by a given amount when it is split. On the other hand, a smaller node further down on the tree may contain relatively dissimilar observations, but because of the node’s small size, a split at this node would only reduce the overall SSE by </w:t></w:r><w:r><w:t xml:space="preserve">so much. </w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="160"/></w:pPr><w:r><w:t>With a fractional penalty, this would not be the case. The large node of similar observations may possess an absolutely larger SSE within the node than the smaller node of dissimilar observations possesses. But the relative reduction in SSE when splitting the smaller node as opposed to the larger node would be greater; thus, the split at the smaller node would be valued more highly than the split at the larger node. This seems appropriate.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="160"/></w:pPr><w:r><w:t>With this in mind, we opted for the fractional penalty.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$r.InsertXML($xml)

Write-Output "Inserted new paragraphs after anchor."
